$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.26"
$ws.Range("E2").Value = "'-2.61%"
$ws.Range("G2").Value = "'16"

$ws.Range("D3").Value = "'37.96"
$ws.Range("E3").Value = "'-4.74%"
$ws.Range("G3").Value = "'16"

$ws.Range("D4").Value = "'5.101"
$ws.Range("E4").Value = "'-0.92%"
$ws.Range("G4").Value = "'16"

$ws.Range("D5").Value = "'0.07864"
$ws.Range("E5").Value = "'-4.22%"
$ws.Range("G5").Value = "'16"

$ws.Range("D6").Value = "'1.960"
$ws.Range("E6").Value = "'-2.71%"
$ws.Range("G6").Value = "'16"

$ws.Range("D7").Value = "'4.365"
$ws.Range("E7").Value = "'1.95%"
$ws.Range("G7").Value = "'16"

$ws.Range("D8").Value = "'8.289"
$ws.Range("E8").Value = "'0.03%"
$ws.Range("G8").Value = "'16"

$ws.Range("D9").Value = "'3.069"
$ws.Range("E9").Value = "'-4.88%"
$ws.Range("G9").Value = "'16"

$ws.Range("D10").Value = "'0.9313"
$ws.Range("E10").Value = "'-0.24%"
$ws.Range("G10").Value = "'16"

$ws.Range("D11").Value = "'0.1369"
$ws.Range("E11").Value = "'-2.44%"
$ws.Range("G11").Value = "'16"

$ws.Range("D12").Value = "'0.1994"
$ws.Range("E12").Value = "'0.87%"
$ws.Range("G12").Value = "'16"

$ws.Range("D13").Value = "'0.08925"
$ws.Range("E13").Value = "'-1.68%"
$ws.Range("G13").Value = "'16"

$ws.Range("E14").Value = "'-3.60%"
$ws.Range("G14").Value = "'16"

$ws.Range("D15").Value = "'0.09725"
$ws.Range("E15").Value = "'-0.87%"
$ws.Range("G15").Value = "'16"

$ws.Range("D16").Value = "'0.001388"
$ws.Range("E16").Value = "'-0.96%"
$ws.Range("G16").Value = "'16"

$ws.Range("D17").Value = "'0.006144"
$ws.Range("E17").Value = "'-4.20%"
$ws.Range("G17").Value = "'16"

$ws.Range("E18").Value = "'1,775.18%"
$ws.Range("G18").Value = "'16"

$ws.Range("D19").Value = "'3.573"
$ws.Range("E19").Value = "'-2.56%"
$ws.Range("G19").Value = "'16"

$ws.Range("D20").Value = "'0.3468"
$ws.Range("E20").Value = "'0.17%"
$ws.Range("G20").Value = "'16"

$ws.Range("D21").Value = "'0.1297"
$ws.Range("E21").Value = "'0.19%"
$ws.Range("G21").Value = "'16"

$ws.Range("D22").Value = "'5.000"
$ws.Range("E22").Value = "'2.04%"
$ws.Range("G22").Value = "'16"

$ws.Range("D23").Value = "'0.2488"
$ws.Range("E23").Value = "'1.42%"
$ws.Range("G23").Value = "'16"

$ws.Range("D24").Value = "'0.04321"
$ws.Range("E24").Value = "'-0.31%"
$ws.Range("G24").Value = "'16"

$ws.Range("D25").Value = "'0.001218"
$ws.Range("E25").Value = "'-0.68%"
$ws.Range("G25").Value = "'16"

$ws.Range("D26").Value = "'0.004564"
$ws.Range("E26").Value = "'-4.18%"
$ws.Range("G26").Value = "'16"

$ws.Range("D27").Value = "'0.0001350"
$ws.Range("E27").Value = "'3.90%"
$ws.Range("G27").Value = "'16"

$ws.Range("G28").Value = "'16"

$ws.Range("G29").Value = "'16"

$ws.Range("G30").Value = "'16"

$ws.Range("G31").Value = "'16"

$ws.Range("G32").Value = "'16"

$ws.Range("G33").Value = "'16"

$ws.Range("G34").Value = "'16"

$ws.Range("G35").Value = "'16"

$ws.Range("G36").Value = "'16"

$ws.Range("G37").Value = "'16"

$ws.Range("G38").Value = "'16"

$ws.Range("D39").Value = "'0.02272"
$ws.Range("E39").Value = "'3.19%"
$ws.Range("G39").Value = "'16"

$ws.Range("D40").Value = "'0.05059"
$ws.Range("E40").Value = "'-3.46%"
$ws.Range("G40").Value = "'16"

$ws.Range("D41").Value = "'0.007482"
$ws.Range("E41").Value = "'-0.46%"
$ws.Range("G41").Value = "'16"

$ws.Range("D42").Value = "'0.009932"
$ws.Range("E42").Value = "'-2.28%"
$ws.Range("G42").Value = "'16"

$ws.Range("D43").Value = "'0.1356"
$ws.Range("E43").Value = "'-1.58%"
$ws.Range("G43").Value = "'16"

$ws.Range("D44").Value = "'0.001980"
$ws.Range("E44").Value = "'-7.85%"
$ws.Range("G44").Value = "'16"

$ws.Range("D45").Value = "'0.008791"
$ws.Range("E45").Value = "'-11.01%"
$ws.Range("G45").Value = "'16"

$ws.Range("E46").Value = "'-0.87%"
$ws.Range("G46").Value = "'16"

$ws.Range("E47").Value = "'-0.11%"
$ws.Range("G47").Value = "'16"

$ws.Range("D48").Value = "'0.003000"
$ws.Range("E48").Value = "'8.34%"
$ws.Range("G48").Value = "'16"

$ws.Range("G49").Value = "'16"

$ws.Range("E50").Value = "'-0.11%"
$ws.Range("G50").Value = "'16"

$ws.Range("E51").Value = "'-0.11%"
$ws.Range("G51").Value = "'16"
